$p = $ppt.ActivePresentation
$old = $p.Slides.Item(17)
Write-Host "old SlideID=$($old.SlideID)"
try {
  $old.SlideID = 999
  Write-Host "set worked, new=$($old.SlideID)"
} catch {
  Write-Host "set failed: $_"
}
